$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("巨力索具", "协鑫集成", "天际股份")
    3  = @("大位科技", "巨力索具", "巨力索具")
    4  = @("掌阅科技", "光线传媒", "协鑫集成")
    5  = @("百川股份", "博纳影业", "博纳影业")
    6  = @("捷成股份", "万向钱潮", "百川股份")
    7  = @("中文在线", "捷成股份", "万向钱潮")
    8  = @("协鑫集成", "格林美", "大位科技")
    9  = @("网宿科技", "掌阅科技", "格林美")
    10 = @("博纳影业", "中文在线", "嘉美包装")
    11 = @("万向钱潮", "大位科技", "光线传媒")
    12 = @("光线传媒", "国际复材", "杭电股份")
    13 = @("格林美", "中钨高新", "特发信息")
    14 = @("国际复材", "网宿科技", "掌阅科技")
    15 = @("海兰信", "再升科技", "欢瑞世纪")
    16 = @("利欧股份", "利欧股份", "横店影视")
    17 = @("横店影视", "百川股份", "浙文互联")
    18 = @("中国巨石", "航天发展", "捷成股份")
    19 = @("特发信息", "海兰信", "网宿科技")
    20 = @("欢瑞世纪", "横店影视", "利欧股份")
    21 = @("天际股份", "中国巨石", "航天发展")
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
